$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Jan_2026"
$ws.Range("E1").Value = "Dec_2025"
$ws.Range("F1").Value = "Nov_2025"

$ws.Range("A2").Value = "INE202B01038"
$ws.Range("B2").Value = "Piramal Finance Ltd"
$ws.Range("D2").Value = 9.539166
$ws.Range("E2").Value = 9.580310000000001
$ws.Range("F2").Value = 8.82103
$ws.Range("G2").Value = -0.04114400000000096
$ws.Range("H2").Value = 0.7181359999999994

$ws.Range("A3").Value = "INE0J1Y01017"
$ws.Range("B3").Value = "Life Insurance Corporation Of India"
$ws.Range("D3").Value = 9.300810999999999
$ws.Range("E3").Value = 8.890637
$ws.Range("F3").Value = 9.086228
$ws.Range("G3").Value = 0.4101739999999996
$ws.Range("H3").Value = 0.2145829999999993

$ws.Range("A4").Value = "INE018A01030"
$ws.Range("B4").Value = "Larsen & Toubro Limited"
$ws.Range("D4").Value = 7.712664
$ws.Range("E4").Value = 7.383763
$ws.Range("F4").Value = 7.18798
$ws.Range("G4").Value = 0.3289010000000001
$ws.Range("H4").Value = 0.5246840000000006

$ws.Range("A5").Value = "INE364U01010"
$ws.Range("B5").Value = "Adani Green Energy Limited"
$ws.Range("D5").Value = 7.546295
$ws.Range("E5").Value = 8.285857999999999
$ws.Range("F5").Value = 3.314505
$ws.Range("G5").Value = -0.7395629999999995
$ws.Range("H5").Value = 4.23179

$ws.Range("A6").Value = "INE758E01017"
$ws.Range("B6").Value = "Jio Financial Services Limited"
$ws.Range("D6").Value = 7.410121
$ws.Range("E6").Value = 7.91723
$ws.Range("F6").Value = 8.028605000000001
$ws.Range("G6").Value = -0.5071089999999998
$ws.Range("H6").Value = -0.6184840000000005

$ws.Range("A7").Value = "INE423A01024"
$ws.Range("B7").Value = "Adani Enterprises Limited"
$ws.Range("D7").Value = 7.163543
$ws.Range("E7").Value = 7.320955
$ws.Range("F7").Value = 2.735111
$ws.Range("G7").Value = -0.1574119999999999
$ws.Range("H7").Value = 4.428432

$ws.Range("A8").Value = "INE814H01029"
$ws.Range("B8").Value = "Adani Power Limited"
$ws.Range("D8").Value = 5.793722
$ws.Range("E8").Value = 5.641521
$ws.Range("F8").Value = 5.684118
$ws.Range("G8").Value = 0.1522009999999998
$ws.Range("H8").Value = 0.109604

$ws.Range("A9").Value = "INE115A01026"
$ws.Range("B9").Value = "LIC Housing Finance Ltd"
$ws.Range("D9").Value = 3.978877
$ws.Range("E9").Value = 3.760505
$ws.Range("F9").Value = 3.738315
$ws.Range("G9").Value = 0.218372
$ws.Range("H9").Value = 0.2405620000000002

$ws.Range("A10").Value = "INE406A01037"
$ws.Range("B10").Value = "Aurobindo Pharma Limited"
$ws.Range("D10").Value = 3.857052
$ws.Range("E10").Value = 3.483121
$ws.Range("F10").Value = 3.528029
$ws.Range("G10").Value = 0.3739309999999998
$ws.Range("H10").Value = 0.3290229999999998

$ws.Range("A11").Value = "INE245A01021"
$ws.Range("B11").Value = "Tata Power Company Limited"
$ws.Range("D11").Value = 2.883411
$ws.Range("E11").Value = 2.754756
$ws.Range("F11").Value = 9.316013
$ws.Range("G11").Value = 0.1286550000000002
$ws.Range("H11").Value = -6.432601999999999

$ws.Range("A12").Value = "INE177F01017"
$ws.Range("B12").Value = "Kovai Medical Center & Hospital Ltd."
$ws.Range("D12").Value = 2.125284
$ws.Range("E12").Value = 2.737368
$ws.Range("F12").Value = 2.95693
$ws.Range("G12").Value = -0.6120839999999999
$ws.Range("H12").Value = -0.8316459999999997

$ws.Range("A13").Value = "INE795G01014"
$ws.Range("B13").Value = "HDFC Life Insurance Co Ltd"
$ws.Range("D13").Value = 2.065748
$ws.Range("E13").Value = 1.953537
$ws.Range("F13").Value = 1.945007
$ws.Range("G13").Value = 0.1122110000000001
$ws.Range("H13").Value = 0.1207410000000002

$ws.Range("A14").Value = "INE424H01027"
$ws.Range("B14").Value = "SUN TV Network Limited"
$ws.Range("D14").Value = 2.052233
$ws.Range("E14").Value = 2.000228
$ws.Range("F14").Value = 1.847762
$ws.Range("G14").Value = 0.0520050000000003
$ws.Range("H14").Value = 0.2044710000000003

$ws.Range("A15").Value = "INE868B01028"
$ws.Range("B15").Value = "NCC Ltd"
$ws.Range("D15").Value = 1.603736
$ws.Range("E15").Value = 1.619954
$ws.Range("F15").Value = 1.690582
$ws.Range("G15").Value = -0.01621799999999984
$ws.Range("H15").Value = -0.08684599999999998

$ws.Range("A16").Value = "INE548A01028"
$ws.Range("B16").Value = "HFCL Limited"
$ws.Range("D16").Value = 1.546887
$ws.Range("E16").Value = 1.409042
$ws.Range("F16").Value = 1.424688
$ws.Range("G16").Value = 0.137845
$ws.Range("H16").Value = 0.1221989999999999

$ws.Range("A17").Value = "INE347A01017"
$ws.Range("B17").Value = "Mangalam Cement Limited"
$ws.Range("D17").Value = 1.44772
$ws.Range("E17").Value = 1.306596
$ws.Range("F17").Value = 1.173441
$ws.Range("G17").Value = 0.1411239999999998
$ws.Range("H17").Value = 0.2742789999999999

$ws.Range("A18").Value = "INE949H01023"
$ws.Range("B18").Value = "Man Infraconstruction Limited"
$ws.Range("D18").Value = 1.230851
$ws.Range("E18").Value = 1.318097
$ws.Range("F18").Value = 1.290651
$ws.Range("G18").Value = -0.08724600000000016
$ws.Range("H18").Value = -0.05980000000000008

$ws.Range("A19").Value = "INE14LE01019"
$ws.Range("B19").Value = "Aditya Birla Lifestyle Brands Limited"
$ws.Range("D19").Value = 1.201118
$ws.Range("E19").Value = 1.345062
$ws.Range("F19").Value = 1.340543
$ws.Range("G19").Value = -0.1439440000000001
$ws.Range("H19").Value = -0.1394250000000001

$ws.Range("A20").Value = "INE0U4701011"
$ws.Range("B20").Value = "Digitide Solutions Limited"
$ws.Range("D20").Value = 1.086307
$ws.Range("E20").Value = 1.153297
$ws.Range("F20").Value = 1.213112
$ws.Range("G20").Value = -0.06699000000000011
$ws.Range("H20").Value = -0.1268050000000001

$ws.Range("A21").Value = "INE095N01031"
$ws.Range("B21").Value = "National Building Construction Corp"
$ws.Range("D21").Value = 0.7053970000000001
$ws.Range("E21").Value = 0.801015
$ws.Range("F21").Value = 0.749797
$ws.Range("G21").Value = -0.09561799999999998
$ws.Range("H21").Value = -0.0444

$ws.Range("A22").Value = "INE942C01045"
$ws.Range("B22").Value = "Gujarat Themis Biosyn Ltd"
$ws.Range("D22").Value = 0.532569
$ws.Range("E22").Value = 0.671127
$ws.Range("F22").Value = 0.605021
$ws.Range("G22").Value = -0.1385580000000001
$ws.Range("H22").Value = -0.07245200000000007

$ws.Range("A23").Value = "INE351A01035"
$ws.Range("B23").Value = "Unichem Laboratories Limited"
$ws.Range("D23").Value = 0.353323
$ws.Range("E23").Value = 0.384464
$ws.Range("F23").Value = 0.397052
$ws.Range("G23").Value = -0.03114099999999997
$ws.Range("H23").Value = -0.04372900000000002

$ws.Range("A24").Value = "INE423A20016"
$ws.Range("B24").Value = "Adani Enterprises Limited Rights"
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0.063677
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = -0.063677

$ws.Range("A25").Value = "INE01TY01017"
$ws.Range("B25").Value = "Canara HSBC Life Insurance Company Ltd"
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0.786067
$ws.Range("F25").Value = 0.637611
$ws.Range("G25").Value = -0.786067
$ws.Range("H25").Value = -0.637611

$ws.Range("A26").Value = "INE062A01020"
$ws.Range("B26").Value = "State Bank of India"
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 0.650393
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = -0.650393

$ws.Range("A27").Value = "INE933K01021"
$ws.Range("B27").Value = "Bajaj Consumer Care Ltd"
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 3.739444
$ws.Range("F27").Value = 3.811694
$ws.Range("G27").Value = -3.739444
$ws.Range("H27").Value = -3.811694
